# Have to do switch
# Insert three new paragraphs at the very start of the document body,
# ahead of the existing (picture) paragraph:
#   1. an empty paragraph (paragraph mark carries <w:noProof/>)
#   2. a paragraph containing the text "Hello, I have to switch"
#   3. another empty paragraph (paragraph mark carries <w:noProof/>)
#
# We build the insertion as a raw WordOpenXML fragment and insert it via
# Range.InsertXML at a range collapsed to the very start of the document,
# so the original picture paragraph is preserved untouched and simply
# pushed down after the new content.

$d = $word.ActiveDocument

$firstPara = $d.Paragraphs.Item(1)
$insertionPoint = $d.Range($firstPara.Range.Start, $firstPara.Range.Start)

$insertXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:pPr><w:rPr><w:noProof/></w:rPr></w:pPr></w:p>
<w:p><w:pPr><w:rPr><w:noProof/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:t>Hello, I have to switch</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:noProof/></w:rPr></w:pPr></w:p>
</w:body></w:document>
</pkg:xmlData></pkg:part></pkg:package>
"@

$insertionPoint.InsertXML($insertXml)
